$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers L1, M1
$ws.Range("L1").Value = "battery_size"
$ws.Range("M1").Value = "total_charged"

# Copy header style (bold, border, centered) from an existing header cell to the new ones
$ws.Range("B1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 2
$ws.Range("E2").Value = 187
$ws.Range("F2").Value = 508.8885531366015
$ws.Range("G2").Value = 206.4294023836775
$ws.Range("H2").Value = 62.66630441947908
$ws.Range("I2").Value = 0.08899340861885814
$ws.Range("J2").Value = 0.1796378076245018
$ws.Range("K2").Value = 0.7313687837566401
$ws.Range("L2").Value = 8372.268247304077
$ws.Range("M2").Value = @'
0        21
1        11
2         0
3        41
4        25
       ... 
1301     47
1302    126
1303     63
1304     32
1305    164
Name: total_charged, Length: 129256, dtype: int64
'@
$ws.Rows.Item(2).AutoFit()

# Row 3
$ws.Range("E3").Value = 188
$ws.Range("F3").Value = 509.7620102496328
$ws.Range("G3").Value = 437.710747419341
$ws.Range("H3").Value = 296.9443778079851
$ws.Range("I3").Value = 0.2470400422971748
$ws.Range("J3").Value = 0.2134461767008358
$ws.Range("K3").Value = 0.5395137810019893
$ws.Range("L3").Value = 8383.164324361531
$ws.Range("M3").Value = @'
0       429
1       123
2       443
3       100
4       684
       ... 
1291    277
1292    224
1293    497
1294     86
1295    303
Name: total_charged, Length: 124197, dtype: int64
'@
$ws.Rows.Item(3).AutoFit()

# Row 4
$ws.Range("E4").Value = 187
$ws.Range("F4").Value = 509.6459523893998
$ws.Range("G4").Value = 392.7007081587785
$ws.Range("H4").Value = 181.43397802114
$ws.Range("I4").Value = 0.1809237634091522
$ws.Range("J4").Value = 0.2159653562384574
$ws.Range("K4").Value = 0.6031108803523904
$ws.Range("L4").Value = 8448.797492483765
$ws.Range("M4").Value = @'
0       182
1       177
2        30
3       129
4        59
       ... 
1216      0
1217    115
1218    119
1219    274
1220     32
Name: total_charged, Length: 124587, dtype: int64
'@
$ws.Rows.Item(4).AutoFit()

# Row 5
$ws.Range("E5").Value = 187
$ws.Range("F5").Value = 509.3009020729499
$ws.Range("G5").Value = 153.1316388413295
$ws.Range("H5").Value = 68.1203251827374
$ws.Range("I5").Value = 0.101577447969993
$ws.Range("J5").Value = 0.1301343401230972
$ws.Range("K5").Value = 0.7682882119069097
$ws.Range("L5").Value = 11432.79114088729
$ws.Range("M5").Value = @'
0        45
1       104
2        32
3         0
4        18
       ... 
1324    100
1325     27
1326     19
1327     55
1328     70
Name: total_charged, Length: 129628, dtype: int64
'@
$ws.Rows.Item(5).AutoFit()

# Row 6
$ws.Range("E6").Value = 187
$ws.Range("F6").Value = 510.0075201937716
$ws.Range("G6").Value = 528.661196971686
$ws.Range("H6").Value = 323.5828757992401
$ws.Range("I6").Value = 0.2488622407138089
$ws.Range("J6").Value = 0.2345833513934074
$ws.Range("K6").Value = 0.5165544078927837
$ws.Range("L6").Value = 11509.45863171522
$ws.Range("M6").Value = @'
0         0
1       126
2       469
3       218
4         0
       ... 
1196    393
1197    121
1198    442
1199    142
1200    141
Name: total_charged, Length: 119741, dtype: int64
'@
$ws.Rows.Item(6).AutoFit()

# Row 7
$ws.Range("E7").Value = 186
$ws.Range("F7").Value = 507.5298715607734
$ws.Range("G7").Value = 453.9725382063149
$ws.Range("H7").Value = 194.5419017719246
$ws.Range("I7").Value = 0.1919068959044306
$ws.Range("J7").Value = 0.2339969984908543
$ws.Range("K7").Value = 0.5740961056047152
$ws.Range("L7").Value = 11459.43465567113
$ws.Range("M7").Value = @'
0       238
1       339
2       113
3       220
4       359
       ... 
1283     73
1284    400
1285    406
1286      0
1287    166
Name: total_charged, Length: 124121, dtype: int64
'@
$ws.Rows.Item(7).AutoFit()

# Row 8
$ws.Range("E8").Value = 188
$ws.Range("F8").Value = 511.0362516374855
$ws.Range("G8").Value = 296.6333166384444
$ws.Range("H8").Value = 56.97803654530541
$ws.Range("I8").Value = 0.08175400538347506
$ws.Range("J8").Value = 0.1591490643505704
$ws.Range("K8").Value = 0.7590969302659545
$ws.Range("L8").Value = 5352.09787381037
$ws.Range("M8").Value = @'
0        89
1        48
2        25
3        71
4        42
       ... 
1254    142
1255     44
1256     55
1257     13
1258     10
Name: total_charged, Length: 126173, dtype: int64
'@
$ws.Rows.Item(8).AutoFit()

# Row 9
$ws.Range("E9").Value = 188
$ws.Range("F9").Value = 510.4176330202803
$ws.Range("G9").Value = 831.4068763008148
$ws.Range("H9").Value = 280.2880539549259
$ws.Range("I9").Value = 0.2077131795030495
$ws.Range("J9").Value = 0.3439404083800179
$ws.Range("K9").Value = 0.4483464121169326
$ws.Range("L9").Value = 5347.592470171842
$ws.Range("M9").Value = @'
0       160
1         0
2       418
3       108
4       699
       ... 
1163     47
1164      0
1165     80
1166    621
1167    379
Name: total_charged, Length: 114345, dtype: int64
'@
$ws.Rows.Item(9).AutoFit()

# Row 10
$ws.Range("E10").Value = 188
$ws.Range("F10").Value = 511.3819334949079
$ws.Range("G10").Value = 693.8127643109007
$ws.Range("H10").Value = 170.5160744111434
$ws.Range("I10").Value = 0.1529007232097687
$ws.Range("J10").Value = 0.3349695881806478
$ws.Range("K10").Value = 0.5121296886095835
$ws.Range("L10").Value = 5358.945395616649
$ws.Range("M10").Value = @'
0         0
1       312
2       325
3        67
4       252
       ... 
1187    193
1188     61
1189    252
1190     50
1191    125
Name: total_charged, Length: 117115, dtype: int64
'@
$ws.Rows.Item(10).AutoFit()
